$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C3").Value = "Y"
$ws.Range("C3").Select()
